$wb = $excel.ActiveWorkbook

# Duplicate Sheet1 (placing the copy after it) so the new sheet inherits
# the same sheetFormatPr/namespace boilerplate Excel originally wrote for
# this workbook, then wipe its contents and repopulate it with the new
# data-driven test table. This mirrors the target: a new worksheet
# "FindNewCarTest" added after Sheet1 (sheetId 2 / rId2), which becomes
# the active tab.
$src = $wb.Worksheets.Item(1)
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "FindNewCarTest"
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = '${browser}'
$ws.Range("B1").Value = '${brandname}'
$ws.Range("C1").Value = '${carheading}'

# Column A (browser) filled first, then column B (brand), then column C
# (heading) - this mirrors the shared-string insertion order recorded in
# the target file (toyota/kia/bmw all before the *Cars headings).
$ws.Range("A2").Value = "chrome"
$ws.Range("A3").Value = "chrome"
$ws.Range("A4").Value = "chrome"

$ws.Range("B2").Value = "toyota"
$ws.Range("B3").Value = "kia"
$ws.Range("B4").Value = "bmw"

$ws.Range("C3").Value = "Kia Cars"
$ws.Range("C2").Value = "Toyta Cars"
$ws.Range("C4").Value = "BMW Cars"

# Column widths for B and C (character-width units through COM).
$ws.Range("B1").ColumnWidth = 11.666666666666666
$ws.Range("C1").ColumnWidth = 16.92

# Leave the selection where the author last left it before saving.
[void]$ws.Range("A5").Select()
